$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-01 18:40:56"

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
